# Update each two-digit multiplication fact (and its product) found in
# the worksheet's table cells. Every "old" string is unique within the
# document, so a plain Find/Replace (MatchWholeWord, no wildcards) is
# used per pair; $ok tracks whether each Find actually located its
# target so a silent miss would be surfaced in the output.

$d = $word.ActiveDocument
$misses = @()

$ok = $d.Content.Find.Execute("38×81=3078", $true, $false, $false, $false, $false, $true, 1, $false, "26×55=1430", 2)
if (-not $ok) { $misses += "38×81=3078" }
$ok = $d.Content.Find.Execute("60×75=4500", $true, $false, $false, $false, $false, $true, 1, $false, "29×40=1160", 2)
if (-not $ok) { $misses += "60×75=4500" }
$ok = $d.Content.Find.Execute("12×34=408", $true, $false, $false, $false, $false, $true, 1, $false, "14×17=238", 2)
if (-not $ok) { $misses += "12×34=408" }
$ok = $d.Content.Find.Execute("90×26=2340", $true, $false, $false, $false, $false, $true, 1, $false, "20×95=1900", 2)
if (-not $ok) { $misses += "90×26=2340" }
$ok = $d.Content.Find.Execute("60×44=2640", $true, $false, $false, $false, $false, $true, 1, $false, "94×56=5264", 2)
if (-not $ok) { $misses += "60×44=2640" }
$ok = $d.Content.Find.Execute("92×73=6716", $true, $false, $false, $false, $false, $true, 1, $false, "70×57=3990", 2)
if (-not $ok) { $misses += "92×73=6716" }
$ok = $d.Content.Find.Execute("52×13=676", $true, $false, $false, $false, $false, $true, 1, $false, "53×82=4346", 2)
if (-not $ok) { $misses += "52×13=676" }
$ok = $d.Content.Find.Execute("30×55=1650", $true, $false, $false, $false, $false, $true, 1, $false, "76×50=3800", 2)
if (-not $ok) { $misses += "30×55=1650" }
$ok = $d.Content.Find.Execute("74×52=3848", $true, $false, $false, $false, $false, $true, 1, $false, "39×68=2652", 2)
if (-not $ok) { $misses += "74×52=3848" }
$ok = $d.Content.Find.Execute("53×95=5035", $true, $false, $false, $false, $false, $true, 1, $false, "32×37=1184", 2)
if (-not $ok) { $misses += "53×95=5035" }
$ok = $d.Content.Find.Execute("51×52=2652", $true, $false, $false, $false, $false, $true, 1, $false, "95×78=7410", 2)
if (-not $ok) { $misses += "51×52=2652" }
$ok = $d.Content.Find.Execute("44×32=1408", $true, $false, $false, $false, $false, $true, 1, $false, "30×97=2910", 2)
if (-not $ok) { $misses += "44×32=1408" }
$ok = $d.Content.Find.Execute("57×11=627", $true, $false, $false, $false, $false, $true, 1, $false, "14×22=308", 2)
if (-not $ok) { $misses += "57×11=627" }
$ok = $d.Content.Find.Execute("99×40=3960", $true, $false, $false, $false, $false, $true, 1, $false, "85×45=3825", 2)
if (-not $ok) { $misses += "99×40=3960" }
$ok = $d.Content.Find.Execute("46×54=2484", $true, $false, $false, $false, $false, $true, 1, $false, "92×40=3680", 2)
if (-not $ok) { $misses += "46×54=2484" }
$ok = $d.Content.Find.Execute("66×33=2178", $true, $false, $false, $false, $false, $true, 1, $false, "44×91=4004", 2)
if (-not $ok) { $misses += "66×33=2178" }
$ok = $d.Content.Find.Execute("26×75=1950", $true, $false, $false, $false, $false, $true, 1, $false, "30×16=480", 2)
if (-not $ok) { $misses += "26×75=1950" }
$ok = $d.Content.Find.Execute("38×72=2736", $true, $false, $false, $false, $false, $true, 1, $false, "32×92=2944", 2)
if (-not $ok) { $misses += "38×72=2736" }
$ok = $d.Content.Find.Execute("77×25=1925", $true, $false, $false, $false, $false, $true, 1, $false, "99×60=5940", 2)
if (-not $ok) { $misses += "77×25=1925" }
$ok = $d.Content.Find.Execute("83×63=5229", $true, $false, $false, $false, $false, $true, 1, $false, "34×23=782", 2)
if (-not $ok) { $misses += "83×63=5229" }
$ok = $d.Content.Find.Execute("78×14=1092", $true, $false, $false, $false, $false, $true, 1, $false, "70×92=6440", 2)
if (-not $ok) { $misses += "78×14=1092" }
$ok = $d.Content.Find.Execute("40×85=3400", $true, $false, $false, $false, $false, $true, 1, $false, "65×37=2405", 2)
if (-not $ok) { $misses += "40×85=3400" }
$ok = $d.Content.Find.Execute("84×18=1512", $true, $false, $false, $false, $false, $true, 1, $false, "23×14=322", 2)
if (-not $ok) { $misses += "84×18=1512" }
$ok = $d.Content.Find.Execute("75×41=3075", $true, $false, $false, $false, $false, $true, 1, $false, "43×67=2881", 2)
if (-not $ok) { $misses += "75×41=3075" }
$ok = $d.Content.Find.Execute("41×77=3157", $true, $false, $false, $false, $false, $true, 1, $false, "41×69=2829", 2)
if (-not $ok) { $misses += "41×77=3157" }

if ($misses.Count -gt 0) {
    Write-Output ("Missed replacements: " + ($misses -join ", "))
} else {
    Write-Output "Replacements applied: 25"
}
